$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the "year" column header (column B) to "season_ending_year"
$ws.Range("B1").Value = "season_ending_year"

# 2) Fill in the previously-blank "birth_year" column (E) with 1985 for every
#    data row (LeBron James' birth year), matching column header already in E1.
for ($r = 2; $r -le 23; $r++) {
    $ws.Cells.Item($r, 5).Value = 1985
}

# 3) Add a new trailing column AM ("calendar_year") holding the calendar year
#    that corresponds to each season's "year" value already present in column B.
$ws.Cells.Item(1, 39).Value = "calendar_year"
# Match the header formatting (bold, centered, bordered) used by the other
# header cells in row 1 by copying it from the adjacent header cell (AL1).
$ws.Cells.Item(1, 38).Copy()
$ws.Cells.Item(1, 39).PasteSpecial(-4122)  # xlPasteFormats

$calendarYears = @(2025, 2024, 2023, 2022, 2021, 2020, 2019, 2018, 2017, 2016, 2015, 2014, 2013, 2012, 2011, 2010, 2009, 2008, 2007, 2006, 2005, 2004)
for ($i = 0; $i -lt $calendarYears.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 39).Value = $calendarYears[$i]
}
